{"js": "// Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Appends \"Briefing:\" to the trailing (empty) bullet under \"Aula 2 \u2013 Estrat\u00e9gias\n// de SEO:\" and then adds the new \"Briefing\" paragraph plus the \"O que\n// aprendemos:\" summary bullets that follow it, mirroring the author's edit.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\n// The last paragraph in the document is the empty (\" \") bullet at level 1\n// right after \"Aula 2 \u2013 Estrat\u00e9gias de SEO:\".\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\n\n// Add the \"Briefing:\" text to the end of that paragraph (keeps the existing\n// \" \" text already present, same as the diff keeps the original run).\nlastParagraph.getRange(\"End\").insertText(\"Briefing:\", \"End\");\n\n// New sub-bullet (level 2) explaining what a briefing is.\nconst briefingParagraph = lastParagraph.insertParagraph(\n  \"Briefing \u00e9 o planejamento em que s\u00e3o definidos, de forma bem detalhada, todos os objetivos do projeto. Portanto, o seu preenchimento \u00e9 indispens\u00e1vel.\",\n  \"After\"\n);\nbriefingParagraph.listItem.level = 2;\n\n// New bullet (level 1) introducing the \"what we learned\" summary.\nconst learnedParagraph = briefingParagraph.insertParagraph(\" O que aprendemos:\", \"After\");\nlearnedParagraph.listItem.level = 1;\n\n// Summary sub-bullets (level 2).\nconst bullet1 = learnedParagraph.insertParagraph(\n  \"O que \u00e9 necess\u00e1rio para iniciar o seu planejamento estrat\u00e9gico de SEO;\",\n  \"After\"\n);\nbullet1.listItem.level = 2;\n\nconst bullet2 = bullet1.insertParagraph(\n  \"O que \u00e9 um briefing, e como fazer um bom briefing para SEO;\",\n  \"After\"\n);\nbullet2.listItem.level = 2;\n\nconst bullet3 = bullet2.insertParagraph(\n  \"O que s\u00e3o palavras-chave, e como iniciar sua pesquisa de palavras para montar o seu conte\u00fado de acordo com o objetivo de marketing.\",\n  \"After\"\n);\nbullet3.listItem.level = 2;\n\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d below.\n#\n# Appends \"Briefing:\" to the trailing (empty) bullet under \"Aula 2 \u2013\n# Estrat\u00e9gias de SEO:\" and then adds the new \"Briefing\" paragraph plus the\n# \"O que aprendemos:\" summary bullets that follow it, mirroring the author's\n# edit.\n\n$d = $word.ActiveDocument\n\n# The last paragraph in the document is the empty (\" \") bullet at list level\n# 2 (ListLevelNumber, 1-based) right after \"Aula 2 \u2013 Estrat\u00e9gias de SEO:\".\n$lastParagraph = $d.Paragraphs.Last\n\n# Add the \"Briefing:\" text to the end of that paragraph (keeps the existing\n# \" \" text already present, same as the diff keeps the original run).\n$lastParagraph.Range.InsertAfter(\"Briefing:\")\n\n# New sub-bullet (ListLevelNumber 3 == ilvl 2) explaining what a briefing is.\n$lastParagraph.Range.InsertParagraphAfter()\n$briefingParagraph = $d.Paragraphs.Last\n$briefingParagraph.Range.Text = \"Briefing \u00e9 o planejamento em que s\u00e3o definidos, de forma bem detalhada, todos os objetivos do projeto. Portanto, o seu preenchimento \u00e9 indispens\u00e1vel.\"\n$briefingParagraph.Range.ListFormat.ListLevelNumber = 3\n\n# New bullet (ListLevelNumber 2 == ilvl 1) introducing the \"what we learned\"\n# summary.\n$briefingParagraph.Range.InsertParagraphAfter()\n$learnedParagraph = $d.Paragraphs.Last\n$learnedParagraph.Range.Text = \" O que aprendemos:\"\n$learnedParagraph.Range.ListFormat.ListLevelNumber = 2\n\n# Summary sub-bullets (ListLevelNumber 3 == ilvl 2).\n$learnedParagraph.Range.InsertParagraphAfter()\n$bullet1 = $d.Paragraphs.Last\n$bullet1.Range.Text = \"O que \u00e9 necess\u00e1rio para iniciar o seu planejamento estrat\u00e9gico de SEO;\"\n$bullet1.Range.ListFormat.ListLevelNumber = 3\n\n$bullet1.Range.InsertParagraphAfter()\n$bullet2 = $d.Paragraphs.Last\n$bullet2.Range.Text = \"O que \u00e9 um briefing, e como fazer um bom briefing para SEO;\"\n$bullet2.Range.ListFormat.ListLevelNumber = 3\n\n$bullet2.Range.InsertParagraphAfter()\n$bullet3 = $d.Paragraphs.Last\n$bullet3.Range.Text = \"O que s\u00e3o palavras-chave, e como iniciar sua pesquisa de palavras para montar o seu conte\u00fado de acordo com o objetivo de marketing.\"\n$bullet3.Range.ListFormat.ListLevelNumber = 3\n\n$d.Save()\n"}
